$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 (shifts old rows 42-43 down to 43-44)
$ws.Rows.Item(42).Insert()

# Fill in the new row 42 with the new weekly record
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 44826
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100102
$ws.Cells.Item(42, 8).Value = "Cítricos"
$ws.Cells.Item(42, 9).Value = 100102006
$ws.Cells.Item(42, 10).Value = "Pomelo"
$ws.Cells.Item(42, 11).Value = "Start Ruby"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 18
$ws.Cells.Item(42, 14).Value = 170000
$ws.Cells.Item(42, 15).Value = 170000
$ws.Cells.Item(42, 16).Value = 170000
$ws.Cells.Item(42, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(42, 18).Value = "Región Metropolitana"
$ws.Cells.Item(42, 19).Value = 486
$ws.Cells.Item(42, 20).Value = 350
